# Add "Neutered Status" answers (column I) for the NCATS-COP01 Transcriptomics
# block of rows (86-145). Every row gets "Yes" except row 104, which gets "No".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($r = 86; $r -le 145; $r++) {
    if ($r -eq 104) {
        $ws.Cells.Item($r, 9).Value = "No"
    } else {
        $ws.Cells.Item($r, 9).Value = "Yes"
    }
}
